# Update the "ex9.1.9 (Linear) M Stationary generator alpha zero" workbook so
# every non-convex experiment except the 5th one is refreshed with new
# generator output values (commit: "expermits todos no convexos menos el
# 5to").
#
# The cells below hold plain numeric-looking text (e.g. "-8.8") that must be
# stored as literal strings (shared-string table entries), NOT as numbers -
# exactly like the rest of the workbook already does. Assigning such a
# string straight to Range.Value lets Excel's normal "smart" input parsing
# kick in and silently reinterprets it as a number (and a bare
# NumberFormat="@" / leading apostrophe trick pulls in a new style record).
# Building the text on a scratch cell via a `="..."` formula, copying it,
# and pasting-values-only onto the destination keeps the string exactly as
# typed (including the '-0.0' style of formatting) without touching
# styles.xml at all.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts the new sheet at tab position 1,
# shifting every existing sheet's (1-based) index by one - so append the
# scratch sheet After the last existing tab instead, to keep all the real
# sheets' positions (and any index-based lookups below) stable.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch = $wb.Worksheets.Add($null, $lastSheet)

function Set-TextValue {
    param($targetSheet, $cellRef, $text)

    # $targetSheet may be a sheet name or a (1-based) tab index.
    $ws = $wb.Worksheets.Item($targetSheet)
    $scratch.Range("A1").Formula = '="' + $text + '"'
    $scratch.Range("A1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Restricciones_del_follower
Set-TextValue "Restricciones_del_follower" "A2" "11.3 - 2x_1 + y_1 - y_2"
Set-TextValue "Restricciones_del_follower" "B2" "-8.8"
Set-TextValue "Restricciones_del_follower" "D2" "0.79"
Set-TextValue "Restricciones_del_follower" "F2" "7.6"

Set-TextValue "Restricciones_del_follower" "A3" "-3.55 + x_1 - 3x_2 + y_2"
Set-TextValue "Restricciones_del_follower" "B3" "1.5499999999999998"
Set-TextValue "Restricciones_del_follower" "D3" "0.09"
Set-TextValue "Restricciones_del_follower" "E3" "5.2"
Set-TextValue "Restricciones_del_follower" "F3" "0"

Set-TextValue "Restricciones_del_follower" "A4" "-9.56 + x_1 + x_2"
Set-TextValue "Restricciones_del_follower" "B4" "6.75"
Set-TextValue "Restricciones_del_follower" "D4" "0.54"
Set-TextValue "Restricciones_del_follower" "E4" "0.8"
Set-TextValue "Restricciones_del_follower" "F4" "0"

# Punto_modificado
Set-TextValue "Punto_modificado" "A2" "6.65"
Set-TextValue "Punto_modificado" "B2" "2.1"
Set-TextValue "Punto_modificado" "C2" "5.2"
Set-TextValue "Punto_modificado" "D2" "3.2"

# Vector_bf
Set-TextValue "Vector_bf" "A2" "3.21"
Set-TextValue "Vector_bf" "A3" "-0.29999999999999993"

# Vector_BF - "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively (same as real Excel),
# so address this one by its (1-based) tab position instead to avoid
# colliding with the "Vector_bf" sheet above.
Set-TextValue 6 "A2" "-4.0"
Set-TextValue 6 "A3" "13.8"
Set-TextValue 6 "A5" "-5.2"

$excel.DisplayAlerts = $false
$scratch.Delete()

# Restore the original active sheet/selection (tab 1 was active before this
# script ran) so the edit doesn't leave stray view-state churn behind.
$wb.Worksheets.Item(1).Activate()
